$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 152.35294
$ws.Range("I39").Value = 54.727272
$ws.Range("J39").Value = 331.33334
$ws.Range("K39").Value = 164.181816
$ws.Range("L39").Value = 994.0000200000001
$ws.Range("M39").Value = 131.818184
$ws.Range("N39").Value = -1586.00002
$ws.Range("H116").Value = 2456.3572
$ws.Range("I116").Value = 2258.9
$ws.Range("K116").Value = 2258.9
$ws.Range("M116").Value = 1183.1
$ws.Range("H132").Value = 4100924.5
$ws.Range("I132").Value = 2349.0176
$ws.Range("J132").Value = 62505624
$ws.Range("K132").Value = 7047.0528
$ws.Range("L132").Value = 187516872
$ws.Range("M132").Value = -4517.0528
$ws.Range("N132").Value = -187521932
$ws.Range("H135").Value = 1173.1389
$ws.Range("I135").Value = 1150.9706
$ws.Range("J135").Value = 1550
$ws.Range("K135").Value = 10358.7354
$ws.Range("L135").Value = 13950
$ws.Range("M135").Value = -7823.735400000001
$ws.Range("N135").Value = -19020
$ws.Range("H137").Value = 967.7353000000001
$ws.Range("I137").Value = 950.5
$ws.Range("J137").Value = 992.3570999999999
$ws.Range("K137").Value = 2851.5
$ws.Range("L137").Value = 2977.0713
$ws.Range("M137").Value = -301.5
$ws.Range("N137").Value = -8077.0713
$ws.Range("H141").Value = 1800.9375
$ws.Range("I141").Value = 1721
$ws.Range("K141").Value = 5163
$ws.Range("M141").Value = 17

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 408.875
$ws.Range("I5").Value = 428.5
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 428.5
$ws.Range("L5").Value = 350
$ws.Range("M5").Value = -316.5
$ws.Range("N5").Value = -574
$ws.Range("H61").Value = 3627.4583
$ws.Range("I61").Value = 3739.9473
$ws.Range("J61").Value = 3200
$ws.Range("K61").Value = 3739.9473
$ws.Range("L61").Value = 3200
$ws.Range("M61").Value = -3527.9473
$ws.Range("N61").Value = -3624
$ws.Range("H74").Value = 2331.6843
$ws.Range("I74").Value = 2433.6667
$ws.Range("J74").Value = 2284.6155
$ws.Range("K74").Value = 2433.6667
$ws.Range("L74").Value = 2284.6155
$ws.Range("M74").Value = -1559.6667
$ws.Range("N74").Value = -4032.6155
$ws.Range("H77").Value = 2331.6843
$ws.Range("I77").Value = 2433.6667
$ws.Range("J77").Value = 2284.6155
$ws.Range("K77").Value = 12168.3335
$ws.Range("L77").Value = 11423.0775
$ws.Range("M77").Value = -7800.333500000001
$ws.Range("N77").Value = -20159.0775
$ws.Range("H136").Value = 3627.4583
$ws.Range("I136").Value = 3739.9473
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 11219.8419
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -8669.841899999999
$ws.Range("N136").Value = -14700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 408.875
$ws.Range("I4").Value = 428.5
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 428.5
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = -313.5
$ws.Range("N4").Value = -580
$ws.Range("H22").Value = 1664.8
$ws.Range("I22").Value = 1527.7778
$ws.Range("J22").Value = 2898
$ws.Range("K22").Value = 1527.7778
$ws.Range("L22").Value = 2898
$ws.Range("M22").Value = -1354.7778
$ws.Range("N22").Value = -3244
$ws.Range("H134").Value = 3675.5454
$ws.Range("I134").Value = 2629.0645
$ws.Range("J134").Value = 5027.25
$ws.Range("K134").Value = 7887.193499999999
$ws.Range("L134").Value = 15081.75
$ws.Range("M134").Value = -5352.193499999999
$ws.Range("N134").Value = -20151.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2843.612
$ws.Range("I31").Value = 1910.5366
$ws.Range("J31").Value = 4315
$ws.Range("K31").Value = 1910.5366
$ws.Range("L31").Value = 4315
$ws.Range("M31").Value = -1615.5366
$ws.Range("N31").Value = -4905
$ws.Range("H34").Value = 2843.612
$ws.Range("I34").Value = 1910.5366
$ws.Range("J34").Value = 4315
$ws.Range("K34").Value = 1910.5366
$ws.Range("L34").Value = 4315
$ws.Range("M34").Value = -1708.5366
$ws.Range("N34").Value = -4719
$ws.Range("H58").Value = 3287.077
$ws.Range("I58").Value = 963.6667
$ws.Range("J58").Value = 5278.5713
$ws.Range("K58").Value = 963.6667
$ws.Range("L58").Value = 5278.5713
$ws.Range("M58").Value = -760.6667
$ws.Range("N58").Value = -5684.5713
$ws.Range("H132").Value = 1799.5814
$ws.Range("I132").Value = 1384.4445
$ws.Range("J132").Value = 3934.5715
$ws.Range("K132").Value = 4153.333500000001
$ws.Range("L132").Value = 11803.7145
$ws.Range("M132").Value = -1623.333500000001
$ws.Range("N132").Value = -16863.7145
$ws.Range("H136").Value = 3287.077
$ws.Range("I136").Value = 963.6667
$ws.Range("J136").Value = 5278.5713
$ws.Range("K136").Value = 2891.0001
$ws.Range("L136").Value = 15835.7139
$ws.Range("M136").Value = -341.0001000000002
$ws.Range("N136").Value = -20935.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3514.1428
$ws.Range("I131").Value = 3484.2144
$ws.Range("J131").Value = 3524.1191
$ws.Range("K131").Value = 10452.6432
$ws.Range("L131").Value = 10572.3573
$ws.Range("M131").Value = -5412.643199999999
$ws.Range("N131").Value = -20652.3573

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 10024.9
$ws.Range("J57").Value = 15570
$ws.Range("L57").Value = 15570
$ws.Range("N57").Value = -17210
$ws.Range("H123").Value = 24968.428
$ws.Range("J123").Value = 24968.428
$ws.Range("L123").Value = 24968.428
$ws.Range("N123").Value = -29868.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2098.4285
$ws.Range("I100").Value = 1546.8572
$ws.Range("J100").Value = 2650
$ws.Range("K100").Value = 1546.8572
$ws.Range("L100").Value = 2650
$ws.Range("M100").Value = -1005.8572
$ws.Range("N100").Value = -3732
$ws.Range("H136").Value = 5977.3335
$ws.Range("I136").Value = 2564.4
$ws.Range("J136").Value = 23042
$ws.Range("K136").Value = 7693.200000000001
$ws.Range("L136").Value = 69126
$ws.Range("M136").Value = -5143.200000000001
$ws.Range("N136").Value = -74226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2646.5
$ws.Range("I132").Value = 2562.1904
$ws.Range("J132").Value = 2807.4546
$ws.Range("K132").Value = 7686.5712
$ws.Range("L132").Value = 8422.363799999999
$ws.Range("M132").Value = -5156.5712
$ws.Range("N132").Value = -13482.3638
$ws.Range("H136").Value = 1795.591
$ws.Range("I136").Value = 956.59375
$ws.Range("J136").Value = 4032.9167
$ws.Range("K136").Value = 2869.78125
$ws.Range("L136").Value = 12098.7501
$ws.Range("M136").Value = -319.78125
$ws.Range("N136").Value = -17198.7501
